$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a literal text value into a cell without letting Excel
# auto-convert number-looking strings (e.g. "1.005") into real numbers.
# Builds a text-producing formula ("=""<text>""") and then converts the
# formula result back to a plain value in place, so the cell keeps its
# original (default) style/number format, matching the source data.
function Set-TextValue {
    param($cell, [string]$text)
    $escaped = $text -replace '"', '""'
    $cell.Formula = '="' + $escaped + '"'
    $cell.Copy()
    $cell.PasteSpecial(-4163)  # xlPasteValues
}

Set-TextValue $ws.Range("D2") '26.561.79'
$ws.Range("E2").Value = '  +6.74%  '
Set-TextValue $ws.Range("D3") '1.743.35'
$ws.Range("E3").Value = '  +5.00%  '
Set-TextValue $ws.Range("D4") '1.005'
$ws.Range("E4").Value = '  +0.32%  '
Set-TextValue $ws.Range("D5") '335.57'
$ws.Range("E5").Value = '  +7.77%  '
Set-TextValue $ws.Range("D6") '1.004'
$ws.Range("E6").Value = '  +0.33%  '
Set-TextValue $ws.Range("D7") '0.3778'
$ws.Range("E7").Value = '  +4.10%  '
Set-TextValue $ws.Range("D8") '48.50'
$ws.Range("E8").Value = '  +2.67%  '
Set-TextValue $ws.Range("D9") '0.3386'
$ws.Range("E9").Value = '  +4.12%  '
$ws.Range("E10").Value = '  +4.73%  '
Set-TextValue $ws.Range("D11") '0.07479'
$ws.Range("E11").Value = '  +5.79%  '
Set-TextValue $ws.Range("D12") '1.004'
$ws.Range("E12").Value = '  +0.51%  '
Set-TextValue $ws.Range("D13") '6.455'
$ws.Range("E13").Value = '  +6.83%  '
$ws.Range("E14").Value = '  +4.55%  '
Set-TextValue $ws.Range("D15") '7.127'
$ws.Range("E15").Value = '  +8.15%  '
Set-TextValue $ws.Range("D16") '1.742.12'
$ws.Range("E16").Value = '  +5.09%  '
Set-TextValue $ws.Range("D17") '0.00001086'
$ws.Range("E17").Value = '  +3.77%  '
Set-TextValue $ws.Range("D18") '0.06713'
$ws.Range("E18").Value = '  +1.72%  '
Set-TextValue $ws.Range("D19") '83.46'
$ws.Range("E19").Value = '  +5.72%  '
$ws.Range("E20").Value = '  +0.22%  '
Set-TextValue $ws.Range("D21") '16.74'
$ws.Range("E21").Value = '  +6.12%  '
Set-TextValue $ws.Range("D22") '6.205'
$ws.Range("E22").Value = '  +5.07%  '
Set-TextValue $ws.Range("D23") '13.03'
$ws.Range("E23").Value = '  +3.98%  '
Set-TextValue $ws.Range("D24") '26.539.97'
$ws.Range("E24").Value = '  +6.87%  '
$ws.Range("E25").Value = '  +0.99%  '
Set-TextValue $ws.Range("D26") '2.462'
$ws.Range("E26").Value = '  +0.46%  '
$ws.Range("B27").Value = 'Monero'
$ws.Range("C27").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextValue $ws.Range("D27") '154.51'
$ws.Range("E27").Value = '  +4.78%  '
$ws.Range("B28").Value = 'ImmutableX'
$ws.Range("C28").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextValue $ws.Range("D28") '1.403'
$ws.Range("E28").Value = '  +16.87%  '
$ws.Range("E29").Value = '  +5.44%  '
Set-TextValue $ws.Range("D30") '1.933.30'
$ws.Range("E30").Value = '  +5.13%  '
Set-TextValue $ws.Range("D31") '132.13'
$ws.Range("E31").Value = '  +5.26%  '
Set-TextValue $ws.Range("D32") '4.136'
$ws.Range("E32").Value = '  +1.57%  '
Set-TextValue $ws.Range("D33") '6.108'
$ws.Range("E33").Value = '  +5.70%  '
Set-TextValue $ws.Range("D34") '0.08668'
$ws.Range("E34").Value = '  +2.66%  '
Set-TextValue $ws.Range("D35") '1.712'
$ws.Range("E35").Value = '  +3.77%  '
Set-TextValue $ws.Range("D36") '12.99'
$ws.Range("E36").Value = '  +5.80%  '
$ws.Range("E37").Value = '  +4.86%  '
Set-TextValue $ws.Range("D38") '0.02358'
$ws.Range("E38").Value = '  +4.43%  '
Set-TextValue $ws.Range("D39") '0.06324'
$ws.Range("E39").Value = '  +4.28%  '
Set-TextValue $ws.Range("D40") '0.2184'
$ws.Range("E40").Value = '  +5.49%  '
Set-TextValue $ws.Range("D41") '8.601'
$ws.Range("E41").Value = '  +3.22%  '
$ws.Range("E42").Value = '  -4.71%  '
Set-TextValue $ws.Range("D43") '0.6236'
$ws.Range("E43").Value = '  +4.85%  '
Set-TextValue $ws.Range("D44") '14.31'
$ws.Range("E44").Value = '  +12.19%  '
Set-TextValue $ws.Range("D45") '1.004'
$ws.Range("E45").Value = '  +0.41%  '
Set-TextValue $ws.Range("D46") '3.932'
$ws.Range("E46").Value = '  +4.28%  '
Set-TextValue $ws.Range("D47") '0.6074'
$ws.Range("E47").Value = '  +7.84%  '
Set-TextValue $ws.Range("D48") '128.56'
$ws.Range("E48").Value = '  +2.65%  '
Set-TextValue $ws.Range("D49") '2.065'
$ws.Range("E49").Value = '  +6.05%  '
Set-TextValue $ws.Range("D50") '0.07253'
$ws.Range("E50").Value = '  +3.84%  '
Set-TextValue $ws.Range("D51") '77.92'
$ws.Range("E51").Value = '  +4.12%  '

$excel.CutCopyMode = 0
